$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText.Replace("1000 Bs = 3.14 = 11919.01 pesos", "1000 Bs = 3.14 = 11871.26 pesos")
$newText = $newText.Replace("11919.01 pesos = 3.12 = 965.31 Bs", "11871.26 pesos = 3.13 = 971.63 Bs")
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10 and N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 317.999
$wsTasas.Range("O10").Value = 3775.05
$wsTasas.Range("N12").Value = 3790
$wsTasas.Range("O12").Value = 310.2
